# Update NatmiData Saa1-Cd36 LR-pair sheet with newly computed TPM-based values.
# (commit: "update scripts wuth new tpm")
#
# The sending-cluster labels for the second and third blocks of rows change
# (MuSCs -> Inflammatory-Mac, Resolving-Mac -> MuSCs) and nearly all of the
# derived statistic columns (E..J, M..T) are recomputed against the new TPM
# matrix. The ligand/receptor symbols (Saa1/Cd36) and the five target-cluster
# labels per block (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac) are
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value2 = [double]"0.2140584197833873"
$ws.Range("J2").Value2 = [double]"0.2140584197833873"
$ws.Range("M2").Value2 = [double]"342.0815733333334"
$ws.Range("N2").Value2 = [double]"1026.24472"
$ws.Range("O2").Value2 = [double]"0.7070495140748019"
$ws.Range("P2").Value2 = [double]"0.7070495140748019"
$ws.Range("Q2").Value2 = [double]"1.710521893857778"
$ws.Range("R2").Value2 = [double]"15.39469704472"
$ws.Range("S2").Value2 = [double]"0.151349901691464"
$ws.Range("T2").Value2 = [double]"0.151349901691464"
$ws.Range("I3").Value2 = [double]"0.2140584197833873"
$ws.Range("J3").Value2 = [double]"0.2140584197833873"
$ws.Range("O3").Value2 = [double]"0.001209381604106811"
$ws.Range("P3").Value2 = [double]"0.001209381604106811"
$ws.Range("S3").Value2 = [double]"0.0002588783150902022"
$ws.Range("T3").Value2 = [double]"0.0002588783150902022"
$ws.Range("I4").Value2 = [double]"0.2140584197833873"
$ws.Range("J4").Value2 = [double]"0.2140584197833873"
$ws.Range("M4").Value2 = [double]"42.10186266666667"
$ws.Range("N4").Value2 = [double]"126.305588"
$ws.Range("O4").Value2 = [double]"0.0870204765782689"
$ws.Range("P4").Value2 = [double]"0.08702047657826889"
$ws.Range("Q4").Value2 = [double]"0.2105233472875556"
$ws.Range("R4").Value2 = [double]"1.894710125588"
$ws.Range("S4").Value2 = [double]"0.01862746570514151"
$ws.Range("T4").Value2 = [double]"0.01862746570514151"
$ws.Range("I5").Value2 = [double]"0.2140584197833873"
$ws.Range("J5").Value2 = [double]"0.2140584197833873"
$ws.Range("M5").Value2 = [double]"3.958736333333333"
$ws.Range("N5").Value2 = [double]"11.876209"
$ws.Range("O5").Value2 = [double]"0.008182324974593572"
$ws.Range("P5").Value2 = [double]"0.008182324974593572"
$ws.Range("Q5").Value2 = [double]"0.01979500124544444"
$ws.Range("R5").Value2 = [double]"0.178155011209"
$ws.Range("S5").Value2 = [double]"0.001751495554215645"
$ws.Range("T5").Value2 = [double]"0.001751495554215645"
$ws.Range("I6").Value2 = [double]"0.2140584197833873"
$ws.Range("J6").Value2 = [double]"0.2140584197833873"
$ws.Range("M6").Value2 = [double]"95.08829366666667"
$ws.Range("N6").Value2 = [double]"285.264881"
$ws.Range("O6").Value2 = [double]"0.1965383027682288"
$ws.Range("P6").Value2 = [double]"0.1965383027682288"
$ws.Range("Q6").Value2 = [double]"0.4754731644312222"
$ws.Range("R6").Value2 = [double]"4.279258479881"
$ws.Range("S6").Value2 = [double]"0.04207067851747599"
$ws.Range("T6").Value2 = [double]"0.04207067851747599"
$ws.Range("A7").Value2 = "Inflammatory-Mac"
$ws.Range("E7").Value2 = [double]"1"
$ws.Range("F7").Value2 = [double]"0.3333333333333333"
$ws.Range("G7").Value2 = [double]"0.000593"
$ws.Range("H7").Value2 = [double]"0.001779"
$ws.Range("I7").Value2 = [double]"0.02538563621056237"
$ws.Range("J7").Value2 = [double]"0.02538563621056237"
$ws.Range("M7").Value2 = [double]"342.0815733333334"
$ws.Range("N7").Value2 = [double]"1026.24472"
$ws.Range("O7").Value2 = [double]"0.7070495140748019"
$ws.Range("P7").Value2 = [double]"0.7070495140748019"
$ws.Range("Q7").Value2 = [double]"0.2028543729866667"
$ws.Range("R7").Value2 = [double]"1.82568935688"
$ws.Range("S7").Value2 = [double]"0.01794890174715782"
$ws.Range("T7").Value2 = [double]"0.01794890174715781"
$ws.Range("A8").Value2 = "Inflammatory-Mac"
$ws.Range("E8").Value2 = [double]"1"
$ws.Range("F8").Value2 = [double]"0.3333333333333333"
$ws.Range("G8").Value2 = [double]"0.000593"
$ws.Range("H8").Value2 = [double]"0.001779"
$ws.Range("I8").Value2 = [double]"0.02538563621056237"
$ws.Range("J8").Value2 = [double]"0.02538563621056237"
$ws.Range("O8").Value2 = [double]"0.001209381604106811"
$ws.Range("P8").Value2 = [double]"0.001209381604106811"
$ws.Range("Q8").Value2 = [double]"0.0003469747763333333"
$ws.Range("R8").Value2 = [double]"0.003122772987"
$ws.Range("S8").Value2 = [double]"3.070092144160188E-05"
$ws.Range("T8").Value2 = [double]"3.070092144160187E-05"
$ws.Range("A9").Value2 = "Inflammatory-Mac"
$ws.Range("E9").Value2 = [double]"1"
$ws.Range("F9").Value2 = [double]"0.3333333333333333"
$ws.Range("G9").Value2 = [double]"0.000593"
$ws.Range("H9").Value2 = [double]"0.001779"
$ws.Range("I9").Value2 = [double]"0.02538563621056237"
$ws.Range("J9").Value2 = [double]"0.02538563621056237"
$ws.Range("M9").Value2 = [double]"42.10186266666667"
$ws.Range("N9").Value2 = [double]"126.305588"
$ws.Range("O9").Value2 = [double]"0.0870204765782689"
$ws.Range("P9").Value2 = [double]"0.08702047657826889"
$ws.Range("Q9").Value2 = [double]"0.02496640456133334"
$ws.Range("R9").Value2 = [double]"0.224697641052"
$ws.Range("S9").Value2 = [double]"0.002209070161285698"
$ws.Range("T9").Value2 = [double]"0.002209070161285697"
$ws.Range("A10").Value2 = "Inflammatory-Mac"
$ws.Range("E10").Value2 = [double]"1"
$ws.Range("F10").Value2 = [double]"0.3333333333333333"
$ws.Range("G10").Value2 = [double]"0.000593"
$ws.Range("H10").Value2 = [double]"0.001779"
$ws.Range("I10").Value2 = [double]"0.02538563621056237"
$ws.Range("J10").Value2 = [double]"0.02538563621056237"
$ws.Range("M10").Value2 = [double]"3.958736333333333"
$ws.Range("N10").Value2 = [double]"11.876209"
$ws.Range("O10").Value2 = [double]"0.008182324974593572"
$ws.Range("P10").Value2 = [double]"0.008182324974593572"
$ws.Range("Q10").Value2 = [double]"0.002347530645666666"
$ws.Range("R10").Value2 = [double]"0.021127775811"
$ws.Range("S10").Value2 = [double]"0.0002077135251616314"
$ws.Range("T10").Value2 = [double]"0.0002077135251616314"
$ws.Range("A11").Value2 = "Inflammatory-Mac"
$ws.Range("E11").Value2 = [double]"1"
$ws.Range("F11").Value2 = [double]"0.3333333333333333"
$ws.Range("G11").Value2 = [double]"0.000593"
$ws.Range("H11").Value2 = [double]"0.001779"
$ws.Range("I11").Value2 = [double]"0.02538563621056237"
$ws.Range("J11").Value2 = [double]"0.02538563621056237"
$ws.Range("M11").Value2 = [double]"95.08829366666667"
$ws.Range("N11").Value2 = [double]"285.264881"
$ws.Range("O11").Value2 = [double]"0.1965383027682288"
$ws.Range("P11").Value2 = [double]"0.1965383027682288"
$ws.Range("Q11").Value2 = [double]"0.05638735814433334"
$ws.Range("R11").Value2 = [double]"0.5074862232989999"
$ws.Range("S11").Value2 = [double]"0.004989249855515619"
$ws.Range("T11").Value2 = [double]"0.004989249855515619"
$ws.Range("A12").Value2 = "MuSCs"
$ws.Range("E12").Value2 = [double]"3"
$ws.Range("F12").Value2 = [double]"1"
$ws.Range("G12").Value2 = [double]"0.01776633333333333"
$ws.Range("H12").Value2 = [double]"0.053299"
$ws.Range("I12").Value2 = [double]"0.7605559440060503"
$ws.Range("J12").Value2 = [double]"0.7605559440060503"
$ws.Range("M12").Value2 = [double]"342.0815733333334"
$ws.Range("N12").Value2 = [double]"1026.24472"
$ws.Range("O12").Value2 = [double]"0.7070495140748019"
$ws.Range("P12").Value2 = [double]"0.7070495140748019"
$ws.Range("Q12").Value2 = [double]"6.077535259031111"
$ws.Range("R12").Value2 = [double]"54.69781733128001"
$ws.Range("S12").Value2 = [double]"0.5377507106361801"
$ws.Range("T12").Value2 = [double]"0.5377507106361801"
$ws.Range("A13").Value2 = "MuSCs"
$ws.Range("E13").Value2 = [double]"3"
$ws.Range("F13").Value2 = [double]"1"
$ws.Range("G13").Value2 = [double]"0.01776633333333333"
$ws.Range("H13").Value2 = [double]"0.053299"
$ws.Range("I13").Value2 = [double]"0.7605559440060503"
$ws.Range("J13").Value2 = [double]"0.7605559440060503"
$ws.Range("O13").Value2 = [double]"0.001209381604106811"
$ws.Range("P13").Value2 = [double]"0.001209381604106811"
$ws.Range("Q13").Value2 = [double]"0.01039539550522222"
$ws.Range("R13").Value2 = [double]"0.093558559547"
$ws.Range("S13").Value2 = [double]"0.0009198023675750074"
$ws.Range("T13").Value2 = [double]"0.0009198023675750074"
$ws.Range("A14").Value2 = "MuSCs"
$ws.Range("E14").Value2 = [double]"3"
$ws.Range("F14").Value2 = [double]"1"
$ws.Range("G14").Value2 = [double]"0.01776633333333333"
$ws.Range("H14").Value2 = [double]"0.053299"
$ws.Range("I14").Value2 = [double]"0.7605559440060503"
$ws.Range("J14").Value2 = [double]"0.7605559440060503"
$ws.Range("M14").Value2 = [double]"42.10186266666667"
$ws.Range("N14").Value2 = [double]"126.305588"
$ws.Range("O14").Value2 = [double]"0.0870204765782689"
$ws.Range("P14").Value2 = [double]"0.08702047657826889"
$ws.Range("Q14").Value2 = [double]"0.7479957260902222"
$ws.Range("R14").Value2 = [double]"6.731961534812"
$ws.Range("S14").Value2 = [double]"0.06618394071184169"
$ws.Range("T14").Value2 = [double]"0.06618394071184169"
$ws.Range("A15").Value2 = "MuSCs"
$ws.Range("E15").Value2 = [double]"3"
$ws.Range("F15").Value2 = [double]"1"
$ws.Range("G15").Value2 = [double]"0.01776633333333333"
$ws.Range("H15").Value2 = [double]"0.053299"
$ws.Range("I15").Value2 = [double]"0.7605559440060503"
$ws.Range("J15").Value2 = [double]"0.7605559440060503"
$ws.Range("M15").Value2 = [double]"3.958736333333333"
$ws.Range("N15").Value2 = [double]"11.876209"
$ws.Range("O15").Value2 = [double]"0.008182324974593572"
$ws.Range("P15").Value2 = [double]"0.008182324974593572"
$ws.Range("Q15").Value2 = [double]"0.07033222927677776"
$ws.Range("R15").Value2 = [double]"0.632990063491"
$ws.Range("S15").Value2 = [double]"0.006223115895216296"
$ws.Range("T15").Value2 = [double]"0.006223115895216296"
$ws.Range("A16").Value2 = "MuSCs"
$ws.Range("E16").Value2 = [double]"3"
$ws.Range("F16").Value2 = [double]"1"
$ws.Range("G16").Value2 = [double]"0.01776633333333333"
$ws.Range("H16").Value2 = [double]"0.053299"
$ws.Range("I16").Value2 = [double]"0.7605559440060503"
$ws.Range("J16").Value2 = [double]"0.7605559440060503"
$ws.Range("M16").Value2 = [double]"95.08829366666667"
$ws.Range("N16").Value2 = [double]"285.264881"
$ws.Range("O16").Value2 = [double]"0.1965383027682288"
$ws.Range("P16").Value2 = [double]"0.1965383027682288"
$ws.Range("Q16").Value2 = [double]"1.689370321379889"
$ws.Range("R16").Value2 = [double]"15.204332892419"
$ws.Range("S16").Value2 = [double]"0.1494783743952372"
$ws.Range("T16").Value2 = [double]"0.1494783743952372"
